$d = $word.ActiveDocument

# Update the date paragraph at the top of the document.
# Direct Range.Text assignment is scoped to the range it is called on
# (unlike Find.Execute with Replace, which this host applies document-wide),
# so it is the safe way to retarget a single run of text while leaving
# the rest of the document (and the run formatting) untouched.
$d.Paragraphs.Item(1).Range.Text = "2024-09-08 Sunday"

# Update the division problems in the table, cell by cell.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "53÷5="  # was "54÷6="
$t.Cell(1, 2).Range.Text = "27÷5="  # was "78÷9="
$t.Cell(1, 3).Range.Text = "18÷7="  # was "64÷7="
$t.Cell(1, 4).Range.Text = "43÷6="  # was "78÷8="
$t.Cell(1, 5).Range.Text = "79÷8="  # was "13÷4="
$t.Cell(5, 1).Range.Text = "83÷8="  # was "22÷8="
$t.Cell(5, 2).Range.Text = "88÷4="  # was "21÷2="
$t.Cell(5, 3).Range.Text = "46÷9="  # was "25÷8="
$t.Cell(5, 4).Range.Text = "78÷7="  # was "12÷8="
$t.Cell(5, 5).Range.Text = "55÷9="  # was "86÷3="
$t.Cell(9, 1).Range.Text = "86÷5="  # was "17÷2="
$t.Cell(9, 2).Range.Text = "89÷5="  # was "83÷8="
$t.Cell(9, 3).Range.Text = "97÷3="  # was "26÷8="
$t.Cell(9, 4).Range.Text = "84÷2="  # was "19÷4="
$t.Cell(9, 5).Range.Text = "78÷7="  # was "89÷2="
$t.Cell(13, 1).Range.Text = "41÷9="  # was "95÷5="
$t.Cell(13, 2).Range.Text = "80÷7="  # was "34÷2="
$t.Cell(13, 3).Range.Text = "96÷5="  # was "22÷7="
$t.Cell(13, 4).Range.Text = "10÷9="  # was "38÷2="
$t.Cell(13, 5).Range.Text = "55÷5="  # was "29÷3="
$t.Cell(17, 1).Range.Text = "13÷7="  # was "91÷7="
$t.Cell(17, 2).Range.Text = "33÷2="  # was "19÷6="
$t.Cell(17, 3).Range.Text = "62÷8="  # was "40÷6="
$t.Cell(17, 4).Range.Text = "83÷5="  # was "93÷2="
$t.Cell(17, 5).Range.Text = "38÷4="  # was "94÷7="

Write-Output "done"
